$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 129 (pushes existing rows 129..249 down to 130..250)
$ws.Rows("129").Insert()

# Populate the new row 129 with the new daily price record.
# Column layout: A Mercado ID, B Mercado, C Region, D Fecha, E Codreg,
# F Categoria ID, G Categoria, H Variedad, I Calidad, J Volumen,
# K Precio minimo, L Precio maximo, M Precio promedio ponderado,
# N Unidad de comercializacion, O Origen, P Precio $/Kg, Q Kg o Unidades,
# R Clasificacion
$ws.Cells.Item(129, 1).Value = 3
$ws.Cells.Item(129, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(129, 3).Value = "Coquimbo"
$ws.Cells.Item(129, 4).Value = 44484
$ws.Cells.Item(129, 5).Value = 5
$ws.Cells.Item(129, 6).Value = 100112017
$ws.Cells.Item(129, 7).Value = "Apio"
$ws.Cells.Item(129, 8).Value = "Americana (o)"
$ws.Cells.Item(129, 9).Value = "Primera"
$ws.Cells.Item(129, 10).Value = 130
$ws.Cells.Item(129, 11).Value = 9000
$ws.Cells.Item(129, 12).Value = 9000
$ws.Cells.Item(129, 13).Value = 9000
$ws.Cells.Item(129, 14).Value = "$/docena de matas"
$ws.Cells.Item(129, 15).Value = "Pan de Azúcar"
$ws.Cells.Item(129, 16).Value = 1500
$ws.Cells.Item(129, 17).Value = 6
$ws.Cells.Item(129, 18).Value = "Hortaliza"
